$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert 2 blank rows above the old row 8 ("Jumlah Penjualan"
#    label) so everything from the old row 8 downwards shifts down by 2
#    (old row 8 -> 10, old row 10 -> 12, ... old row 18 -> 20). Excel's
#    Insert() copies the formatting of the row above into the freshly
#    inserted rows, so clear them straight back to "no formatting" - the
#    edited template has nothing in rows 7-9/11.
# ---------------------------------------------------------------------------
$ws.Rows("7:8").Insert()
$ws.Rows("7:8").Clear()

# ---------------------------------------------------------------------------
# 2) Drop the now-empty helper cells in column C that the cleaned-up
#    template no longer carries (they were blank, formatting-only cells).
# ---------------------------------------------------------------------------
$ws.Range("C13").Clear()
$ws.Range("C14").Clear()
$ws.Range("C16").Clear()
$ws.Range("C18").Clear()
$ws.Range("C20").Clear()

# ---------------------------------------------------------------------------
# 3) Match the editor's final cursor position/selection.
# ---------------------------------------------------------------------------
$ws.Range("D14").Select()
